# Updates for version 25.2409.5030 - "Work" scratch sheet additions
# Bug fixes (cut edit form start position/coordinate) + improvements
# (3D preview of the rendered path, default material thickness).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work")
$ws.Activate()

# Row 78 label written first (matches original authoring order / shared-string table)
$ws.Range("A78").Value = "Bottom:"

# Row 74: Dimension: / value
$ws.Range("A74").Value = "Dimension:"
$ws.Range("B74").Value = 100

# Row 76: column headers for the Top/Bottom/Center table
$ws.Range("C76").Value = "Bottom/Up"
$ws.Range("D76").Value = "Center/Dn"
$ws.Range("E76").Value = "Center/Up"
$ws.Range("B76").Value = "Bottom/Dn"
$ws.Range("F76").Value = "Top/Dn"
$ws.Range("G76").Value = "Top/Up"

# Row 77: Top:
$ws.Range("A77").Value = "Top:"
$ws.Range("B77").Formula = "=0-B74"
$ws.Range("C77").Formula = "=B74"
$ws.Range("D77").Formula = "=0-(B74/2)"
$ws.Range("E77").Formula = "=B74/2"
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0

# Row 78: Bottom: (remaining cells; label already set above)
$ws.Range("B78").Value = 0
$ws.Range("C78").Value = 0
$ws.Range("D78").Formula = "=B74/2"
$ws.Range("E78").Formula = "=0-(B74/2)"
$ws.Range("F78").Formula = "=B74"
$ws.Range("G78").Formula = "=0-B74"

# Row 80/81: Drawn top / bottom
$ws.Range("A80").Value = "Drw Top:"
$ws.Range("B80").Formula = "=B74"
$ws.Range("A81").Value = "Drw Bot:"
$ws.Range("B81").Value = 0

# Row 83-91: Convert Range helper block
$ws.Range("A84").Value = "Value"
$ws.Range("B84").Value = 12
$ws.Range("A83").Value = "Convert Range"
$ws.Range("A85").Value = "fromMin"
$ws.Range("B85").Value = 50
$ws.Range("A87").Value = "toMin"
$ws.Range("B87").Value = 0
$ws.Range("A86").Value = "fromMax"
$ws.Range("B86").Value = -50
$ws.Range("A88").Value = "toMax"
$ws.Range("B88").Value = 100
$ws.Range("A91").Value = "Result"
$ws.Range("B91").Formula = "=(((B84-B85)*B90)/B89)+B87"
$ws.Range("A89").Value = "fromRange"
$ws.Range("B89").Formula = "=B86-B85"
$ws.Range("A90").Value = "toRange"
$ws.Range("B90").Formula = "=B88-B87"

# Row 94-101: Layer / Start / End / Depth / Type table (3D preview data)
$ws.Range("A94").Value = "Layer"
$ws.Range("C94").Value = "Start"
$ws.Range("F94").Value = "End"
$ws.Range("I94").Value = "Depth"
$ws.Range("K94").Value = "Type"

$ws.Range("A95").Value = 0
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 0
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 317.5
$ws.Range("I95").Value = 0
$ws.Range("K95").Value = "Transit"

$ws.Range("A96").Value = 0
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 317.5
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 139.7
$ws.Range("I96").Value = 1.5875
$ws.Range("K96").Value = "Plot"

$ws.Range("A97").Value = 0
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 139.7
$ws.Range("F97").Value = 165.1
$ws.Range("G97").Value = 139.7
$ws.Range("I97").Value = 1.5875
$ws.Range("K97").Value = "Plot"

$ws.Range("A98").Value = 0
$ws.Range("C98").Value = 165.1
$ws.Range("D98").Value = 139.7
$ws.Range("F98").Value = 165.1
$ws.Range("G98").Value = -317.5
$ws.Range("I98").Value = 1.5875
$ws.Range("K98").Value = "Plot"

$ws.Range("A99").Value = 1
$ws.Range("C99").Value = 165.1
$ws.Range("D99").Value = -317.5
$ws.Range("F99").Value = 165.1
$ws.Range("G99").Value = 139.7
$ws.Range("I99").Value = 3.175
$ws.Range("K99").Value = "Plot"

$ws.Range("A100").Value = 1
$ws.Range("C100").Value = 165.1
$ws.Range("D100").Value = 139.7
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 139.7
$ws.Range("I100").Value = 3.175
$ws.Range("K100").Value = "Plot"

$ws.Range("A101").Value = 1
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 139.7
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 317.5
$ws.Range("I101").Value = 3.175
$ws.Range("K101").Value = "Plot"

# Row 103-106: Segments / PointA / PointB / Type table (3D preview path)
$ws.Range("A103").Value = "Segments"
$ws.Range("B103").Value = "Point A"
$ws.Range("E103").Value = "PointB"
$ws.Range("H103").Value = "Type"

$ws.Range("B104").Value = "{0.000,0.000,24.994}"
$ws.Range("E104").Value = "{0.000,317.500,24.994}"
$ws.Range("H104").Value = "Transit"

$ws.Range("B105").Value = "{0.000,317.500,24.994}"
$ws.Range("E105").Value = "{0.000,317.500,3.175}"
$ws.Range("H105").Value = "Transit"

$ws.Range("B106").Value = "{0.000,317.500,3.175}"
$ws.Range("E106").Value = "{0.000,317.500,1.588}"
$ws.Range("H106").Value = "Plunge"


# ---- Update the view: frozen-pane scroll position + active selection ----
$ws.Range("A94").Select()
$excel.ActiveWindow.ScrollRow = 94
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B107").Select()
